$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.778.59"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").Value = "1.627.51"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5057"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2573"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.47%  "
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.37"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07790"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.254"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.00%  "
$ws.Range("D13").Value = "1.628.58"
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("D14").Value = "1.852.06"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5569"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.12%  "
$ws.Range("D17").Value = "0.0₅7537"
$ws.Range("E17").Value = "  -2.69%  "
$ws.Range("D18").Value = "25.802.34"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.66%  "
$ws.Range("E21").Value = "  -3.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.806"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.991"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.33%  "
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.778"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1264"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.727"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.51%  "
$ws.Range("E30").Value = "  -0.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04866"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.277"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.24%  "
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.556"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.8935"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.560"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("D38").Value = "1.128.06"
$ws.Range("E38").Value = "  +2.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5461"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01559"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.001"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.560"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7963"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.50%  "
$ws.Range("D45").Value = "1.781.07"
$ws.Range("D46").Value = "0.0₈112"
$ws.Range("E46").Value = "  -7.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4438"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.21"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05056"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.615"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.000"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.18%  "
